# Scheduled runner update: refresh market-price-derived profit figures
# (currentAveragePrice / NQ / HQ price & profit columns) across the
# per-job leve-profit sheets. Values below come from the upstream
# market-data refresh; no formulas/formatting change.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 67188.125
$ws.Range("I62").Value = 146672
$ws.Range("J62").Value = 5367.3335
$ws.Range("K62").Value = 146672
$ws.Range("L62").Value = 5367.3335
$ws.Range("M62").Value = -146048
$ws.Range("N62").Value = -6615.3335
$ws.Range("H64").Value = 5896.923
$ws.Range("I64").Value = 8960
$ws.Range("K64").Value = 8960
$ws.Range("M64").Value = -8712
$ws.Range("H65").Value = 67188.125
$ws.Range("I65").Value = 146672
$ws.Range("J65").Value = 5367.3335
$ws.Range("K65").Value = 733360
$ws.Range("L65").Value = 26836.6675
$ws.Range("M65").Value = -730240
$ws.Range("N65").Value = -33076.6675
$ws.Range("H67").Value = 5896.923
$ws.Range("I67").Value = 8960
$ws.Range("K67").Value = 8960
$ws.Range("M67").Value = -8102
$ws.Range("H92").Value = 6334.6665
$ws.Range("I92").Value = 6334.6665
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 6334.6665
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -5086.6665
$ws.Range("H107").Value = 89.75
$ws.Range("I107").Value = 89.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 89.75
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 1830.25
$ws.Range("H138").Value = 3559.22
$ws.Range("I138").Value = 694.6061
$ws.Range("J138").Value = 4970.1494
$ws.Range("K138").Value = 2083.8183
$ws.Range("L138").Value = 14910.4482
$ws.Range("M138").Value = 3056.1817
$ws.Range("N138").Value = -25190.4482

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1044.8948
$ws.Range("I74").Value = 1058.5
$ws.Range("K74").Value = 1058.5
$ws.Range("M74").Value = -184.5
$ws.Range("H77").Value = 1044.8948
$ws.Range("I77").Value = 1058.5
$ws.Range("K77").Value = 5292.5
$ws.Range("M77").Value = -924.5
$ws.Range("H132").Value = 3669.8572
$ws.Range("I132").Value = 2673
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 8019
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -5489
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2199.75
$ws.Range("I86").Value = 2020.7368
$ws.Range("J86").Value = 2880
$ws.Range("K86").Value = 2020.7368
$ws.Range("L86").Value = 2880
$ws.Range("M86").Value = -897.7367999999999
$ws.Range("N86").Value = -5126
$ws.Range("H89").Value = 2199.75
$ws.Range("I89").Value = 2020.7368
$ws.Range("J89").Value = 2880
$ws.Range("K89").Value = 10103.684
$ws.Range("L89").Value = 14400
$ws.Range("M89").Value = -4487.683999999999
$ws.Range("N89").Value = -25632
$ws.Range("H105").Value = 5537.273
$ws.Range("I105").Value = 6156.6665
$ws.Range("J105").Value = 2750
$ws.Range("K105").Value = 6156.6665
$ws.Range("L105").Value = 2750
$ws.Range("M105").Value = -4409.6665
$ws.Range("N105").Value = -6244
$ws.Range("H134").Value = 95008.27
$ws.Range("I134").Value = 4389.143
$ws.Range("J134").Value = 253591.75
$ws.Range("K134").Value = 13167.429
$ws.Range("L134").Value = 760775.25
$ws.Range("M134").Value = -10632.429
$ws.Range("N134").Value = -765845.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2350.6667
$ws.Range("I31").Value = 2425.9812
$ws.Range("J31").Value = 2140.5789
$ws.Range("K31").Value = 2425.9812
$ws.Range("L31").Value = 2140.5789
$ws.Range("M31").Value = -2130.9812
$ws.Range("N31").Value = -2730.5789
$ws.Range("H34").Value = 2350.6667
$ws.Range("I34").Value = 2425.9812
$ws.Range("J34").Value = 2140.5789
$ws.Range("K34").Value = 2425.9812
$ws.Range("L34").Value = 2140.5789
$ws.Range("M34").Value = -2223.9812
$ws.Range("N34").Value = -2544.5789
$ws.Range("H132").Value = 2318.5
$ws.Range("I132").Value = 1226.25
$ws.Range("K132").Value = 3678.75
$ws.Range("M132").Value = -1148.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 49300.477
$ws.Range("I139").Value = 56683.89
$ws.Range("K139").Value = 170051.67
$ws.Range("M139").Value = -164911.67

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4602.8696
$ws.Range("I70").Value = 4234
$ws.Range("J70").Value = 4886.615
$ws.Range("K70").Value = 4234
$ws.Range("L70").Value = 4886.615
$ws.Range("M70").Value = -3964
$ws.Range("N70").Value = -5426.615
$ws.Range("H73").Value = 4602.8696
$ws.Range("I73").Value = 4234
$ws.Range("J73").Value = 4886.615
$ws.Range("K73").Value = 4234
$ws.Range("L73").Value = 4886.615
$ws.Range("M73").Value = -3298
$ws.Range("N73").Value = -6758.615
$ws.Range("H126").Value = 1819.6
$ws.Range("I126").Value = 1843
$ws.Range("J126").Value = 1772.8
$ws.Range("K126").Value = 5529
$ws.Range("L126").Value = 5318.4
$ws.Range("M126").Value = -3059
$ws.Range("N126").Value = -10258.4
$ws.Range("H132").Value = 8931784
$ws.Range("I132").Value = 17861090
$ws.Range("J132").Value = 2477.8572
$ws.Range("K132").Value = 53583270
$ws.Range("L132").Value = 7433.571599999999
$ws.Range("M132").Value = -53580740
$ws.Range("N132").Value = -12493.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 591.8095
$ws.Range("I22").Value = 454.46667
$ws.Range("J22").Value = 935.1667
$ws.Range("K22").Value = 454.46667
$ws.Range("L22").Value = 935.1667
$ws.Range("M22").Value = -159.46667
$ws.Range("N22").Value = -1525.1667
$ws.Range("H27").Value = 591.8095
$ws.Range("I27").Value = 454.46667
$ws.Range("J27").Value = 935.1667
$ws.Range("K27").Value = 454.46667
$ws.Range("L27").Value = 935.1667
$ws.Range("M27").Value = -347.46667
$ws.Range("N27").Value = -1149.1667
$ws.Range("H82").Value = 3150.3333
$ws.Range("I82").Value = 3380.4
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 3380.4
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -3019.4
$ws.Range("N82").Value = -2722
$ws.Range("H85").Value = 3150.3333
$ws.Range("I85").Value = 3380.4
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 3380.4
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = -2132.4
$ws.Range("N85").Value = -4496
$ws.Range("H132").Value = 2090
$ws.Range("I132").Value = 1886.9143
$ws.Range("J132").Value = 2736.182
$ws.Range("K132").Value = 5660.742899999999
$ws.Range("L132").Value = 8208.545999999998
$ws.Range("M132").Value = -3130.742899999999
$ws.Range("N132").Value = -13268.546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 7097.5
$ws.Range("I45").Value = 4569
$ws.Range("J45").Value = 9626
$ws.Range("K45").Value = 4569
$ws.Range("L45").Value = 9626
$ws.Range("M45").Value = -4078
$ws.Range("N45").Value = -10608
$ws.Range("H136").Value = 1239.2
$ws.Range("I136").Value = 1124.4478
$ws.Range("J136").Value = 1830.6154
$ws.Range("K136").Value = 3373.3434
$ws.Range("L136").Value = 5491.8462
$ws.Range("M136").Value = -823.3433999999997
$ws.Range("N136").Value = -10591.8462
